# Reorders the data rows (2-24) of the active sheet according to a fixed
# permutation, while leaving the columns that are identical across every
# row (A, B, C, E, F, G, H, I, J) untouched. Columns D and K..T are the
# ones that actually vary row to row and therefore are the ones that need
# to be physically rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 24

# Columns (by letter) whose values differ row-to-row and therefore need to
# be copied across as part of the reorder.
$cols = @("D", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# target row -> source row (1-based worksheet row numbers)
$mapping = @{
    2  = 11
    3  = 8
    4  = 14
    5  = 7
    6  = 19
    7  = 23
    8  = 3
    9  = 20
    10 = 21
    11 = 22
    12 = 4
    13 = 9
    14 = 10
    15 = 24
    16 = 17
    17 = 18
    18 = 16
    19 = 2
    20 = 6
    21 = 12
    22 = 13
    23 = 15
    24 = 5
}

# Snapshot the original values for the columns that will move, keyed by
# source row number, before any writes happen. Value2 is used for reads
# because it reliably returns the underlying scalar (numeric serials for
# dates, plain strings for text) in this runtime.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Write the snapshotted values back into their new (target) rows.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $mapping[$r]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $srcVals[$c]
    }
}
